# Updates the "北京-漫展信息" workbook to the scrape taken at commit 456a3b4.
#   Sheet 1 "展览"     : a handful of "想去人数" (F) refreshes, plus one new
#                        event (北京·排球少年同人ONLY) inserted as row 33,
#                        pushing the previous rows 33-37 down to 34-38.
#   Sheet 2 "演出"     : "想去人数" (F) refreshes only.
#   Sheet 3 "本地生活" : "想去人数" (F) refreshes only.
#   Sheet 4 "全部类型" : "想去人数" (F) refreshes only (mirrors sheets 1-3).

$wb = $excel.ActiveWorkbook

function Set-NumberCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# Writes $val as literal text, even when it looks like a date (e.g.
# "2024-12-07") or a big number-ish string - matches the workbook's existing
# convention of storing every column as plain text except F/G.
function Set-TextCell($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
}

function Set-EventRow($ws, $row, $idx, $start, $name, $place, $timeRange, $want, $price, $link, $cover) {
    Set-NumberCell $ws $row 1 $idx
    Set-TextCell   $ws $row 2 $start
    Set-TextCell   $ws $row 3 $name
    Set-TextCell   $ws $row 4 $place
    Set-TextCell   $ws $row 5 $timeRange
    Set-NumberCell $ws $row 6 $want
    Set-NumberCell $ws $row 7 $price
    Set-TextCell   $ws $row 8 $link
    Set-TextCell   $ws $row 9 $cover
}

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Plain "想去人数" refreshes (row numbers unaffected by the later insert,
# since they're all above row 33).
Set-NumberCell $ws1 4  6 5894
Set-NumberCell $ws1 10 6 18
Set-NumberCell $ws1 12 6 675
Set-NumberCell $ws1 13 6 1584
Set-NumberCell $ws1 14 6 1584
Set-NumberCell $ws1 15 6 1556
Set-NumberCell $ws1 16 6 551
Set-NumberCell $ws1 17 6 149
Set-NumberCell $ws1 18 6 617
Set-NumberCell $ws1 19 6 4434
Set-NumberCell $ws1 20 6 31
Set-NumberCell $ws1 22 6 3336
Set-NumberCell $ws1 24 6 10
Set-NumberCell $ws1 25 6 47
Set-NumberCell $ws1 26 6 2306
Set-NumberCell $ws1 28 6 334
Set-NumberCell $ws1 30 6 452
Set-NumberCell $ws1 31 6 1225

# New event: insert a row at 33, pushing the former rows 33-37 to 34-38.
$ws1.Rows.Item(33).Insert()

# The insert synthesizes a stray style for the new A33 - restore it to match
# the bold/centered/bordered "index" style used by every other A-column cell
# by pulling the format from the row right below (itself just shifted down
# from the old row 33, so it still carries the original style).
$ws1.Cells.Item(34, 1).Copy()
$ws1.Cells.Item(33, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

Set-EventRow $ws1 33 32 "2024-12-07" "北京·排球少年同人ONLY" `
    "永外高庄138号 北京大红门国际会展中心" "2024.12.07 10:00-12.07 17:00" `
    0 60 "https://show.bilibili.com/platform/detail.html?id=93501" `
    "//i1.hdslb.com/bfs/openplatform/202410/Nq2SuSWE1728971583727.jpeg"

# Rows 34-38 keep the same event data they had as rows 33-37, but with
# the "想去人数" (F) counters bumped to the newly-scraped totals.
Set-EventRow $ws1 34 33 "2024-12-14" "北京·thebONE×GOJO超次元动漫游戏嘉年华17th" `
    "B1层西区会员活动中心（阳坊涮肉对面） 北投购物公园" "2024.12.14 09:30-12.15 17:00" `
    1 1 "https://show.bilibili.com/platform/detail.html?id=93495" `
    "//i0.hdslb.com/bfs/openplatform/202410/Fzz24Usj1728969298701.jpeg"

Set-EventRow $ws1 35 34 "2024-12-14" "北京·奇想派对第五届" `
    "学清路38号金码大厦B座(六道口地铁站B东北口步行110米) BOM嘻番里" "2024.12.14 10:00-12.15 17:30" `
    51 45 "https://show.bilibili.com/platform/detail.html?id=91077" `
    "//i1.hdslb.com/bfs/openplatform/202408/zMayUoC81724229782742.jpeg"

Set-EventRow $ws1 36 35 "2024-12-28" "北京·第20届IJOY漫展xCGF游戏节" `
    "天辰东路7号 北京国家会议中心" "2024.12.28 09:00-12.29 17:00" `
    1208 8.800000000000001 "https://show.bilibili.com/platform/detail.html?id=92633" `
    "//i0.hdslb.com/bfs/openplatform/202409/EQg8HwjJ1726734597607.jpeg"

Set-EventRow $ws1 37 36 "2025-01-17" " 北京·第21届IJOY漫展xCGF游戏节" `
    "天辰东路7号 北京国家会议中心" "2025.01.17 09:00-01.19 17:00" `
    1192 8.800000000000001 "https://show.bilibili.com/platform/detail.html?id=92634" `
    "//i0.hdslb.com/bfs/openplatform/202409/ASXIizNW1726735204415.jpeg"

Set-EventRow $ws1 38 37 "2025-04-19" "北京·可行中国动漫游戏节" `
    "焦化路甲18号 东进国际中心" "2025.04.19 09:00-04.20 18:00" `
    81 85 "https://show.bilibili.com/platform/detail.html?id=92495" `
    "//i1.hdslb.com/bfs/openplatform/202409/28QBTqAo1726293348310.jpeg"

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performances) - "想去人数" refreshes only.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
Set-NumberCell $ws2 10 6 92
Set-NumberCell $ws2 11 6 13
Set-NumberCell $ws2 12 6 106
Set-NumberCell $ws2 15 6 42
Set-NumberCell $ws2 19 6 299

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local Life) - "想去人数" refreshes only.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
Set-NumberCell $ws3 3 6 671
Set-NumberCell $ws3 4 6 181
Set-NumberCell $ws3 5 6 269

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) - mirrors the "想去人数" refreshes above.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
Set-NumberCell $ws4 6  6 671
Set-NumberCell $ws4 7  6 181
Set-NumberCell $ws4 8  6 5894
Set-NumberCell $ws4 18 6 92
Set-NumberCell $ws4 20 6 13
Set-NumberCell $ws4 21 6 18
Set-NumberCell $ws4 23 6 1584
Set-NumberCell $ws4 24 6 106
Set-NumberCell $ws4 25 6 1556
Set-NumberCell $ws4 26 6 551
Set-NumberCell $ws4 27 6 149
Set-NumberCell $ws4 28 6 617
Set-NumberCell $ws4 29 6 4434
Set-NumberCell $ws4 31 6 3336
Set-NumberCell $ws4 33 6 47
Set-NumberCell $ws4 35 6 2306
Set-NumberCell $ws4 37 6 334
Set-NumberCell $ws4 39 6 452
Set-NumberCell $ws4 40 6 1225
Set-NumberCell $ws4 42 6 299
Set-NumberCell $ws4 47 6 1208
Set-NumberCell $ws4 49 6 1192
Set-NumberCell $ws4 50 6 81

Write-Output "edit complete"
